{"js": "// Fixed #418 Empty AQL expressions generate empty lines.\n// The document contains an empty paragraph (leftover from an empty AQL\n// expression) whose run is styled with an orange font color. That whole\n// paragraph must be removed so the surrounding text flows directly from\n// \"Checks ImageServices registration :\" to \"End of demonstration.\".\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text,font/color\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const paragraph = paragraphs.items[i];\n  const isEmpty = paragraph.text === \"\";\n  const color = (paragraph.font.color || \"\").toUpperCase();\n  if (isEmpty && color === \"#E36C0A\") {\n    paragraph.delete();\n    await context.sync();\n    break;\n  }\n}\n", "ps1": "# Fixed #418 Empty AQL expressions generate empty lines.\n# The document has a leftover empty paragraph (from an empty AQL\n# expression) whose run carries an explicit (non-automatic) orange font\n# color. Remove that whole paragraph so the text flows directly from\n# \"Checks ImageServices registration :\" to \"End of demonstration.\".\n\n$d = $word.ActiveDocument\n\n$wdColorAutomatic = -16777216\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    $text = $p.Range.Text\n    $trimmed = $text.TrimEnd([char]13, [char]7)\n    if ($trimmed -eq \"\" -and $p.Range.Font.Color -ne $wdColorAutomatic) {\n        $target = $p\n        break\n    }\n}\n\nif ($target -ne $null) {\n    $target.Range.Delete()\n}\n"}
